# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.972.83"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "1.893.04"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").Value = "'0.7223"
$ws.Range("E5").Value = "  -6.70%  "

$ws.Range("D6").Value = "'241.89"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.3082"
$ws.Range("E8").Value = "  -3.91%  "

$ws.Range("D9").Value = "'25.94"
$ws.Range("E9").Value = "  -6.98%  "

$ws.Range("D10").Value = "'0.06866"
$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07938"
$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7650"
$ws.Range("E12").Value = "  -2.25%  "

$ws.Range("D13").Value = "1.884.84"
$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").Value = "'5.230"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").Value = "'90.64"
$ws.Range("E15").Value = "  -4.40%  "

$ws.Range("D16").Value = "29.959.58"
$ws.Range("E16").Value = "  -1.32%  "

$ws.Range("D17").Value = "'14.06"
$ws.Range("E17").Value = "  -3.19%  "

$ws.Range("D18").Value = "'5.725"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "'0.000007712"
$ws.Range("E19").Value = "  -3.35%  "

$ws.Range("D20").Value = "'237.56"
$ws.Range("E20").Value = "  -7.04%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.133.44"
$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "'6.835"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").Value = "'9.256"
$ws.Range("E25").Value = "  -3.27%  "

$ws.Range("D26").Value = "'165.07"
$ws.Range("E26").Value = "  +0.87%  "

$ws.Range("D27").Value = "'18.87"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").Value = "'0.1265"
$ws.Range("E28").Value = "  -6.75%  "

$ws.Range("D29").Value = "'2.005"
$ws.Range("E29").Value = "  -12.22%  "

$ws.Range("D30").Value = "'1.360"
$ws.Range("E30").Value = "  -0.86%  "

$ws.Range("D31").Value = "'1.530"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").Value = "'4.278"
$ws.Range("E32").Value = "  -3.34%  "

$ws.Range("D33").Value = "'4.046"
$ws.Range("E33").Value = "  -2.11%  "

$ws.Range("D34").Value = "'0.05063"
$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("D35").Value = "'1.264"
$ws.Range("E35").Value = "  -1.85%  "

$ws.Range("D36").Value = "'0.7310"
$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("D37").Value = "'2.736"
$ws.Range("E37").Value = "  -1.37%  "

$ws.Range("D38").Value = "'0.01906"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").Value = "'2.767"
$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("D40").Value = "'6.306"
$ws.Range("E40").Value = "  -2.11%  "

$ws.Range("D41").Value = "'74.27"
$ws.Range("E41").Value = "  -5.74%  "

$ws.Range("D42").Value = "'0.4395"
$ws.Range("E42").Value = "  -2.54%  "

$ws.Range("D43").Value = "'1.909"
$ws.Range("E43").Value = "  -3.36%  "

$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value = "'0.8361"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'100.51"
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("D47").Value = "'7.527"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "'9.766"
$ws.Range("E48").Value = "  -0.48%  "

$ws.Range("D49").Value = "'37.30"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").Value = "2.054.17"
$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "'932.42"
$ws.Range("E51").Value = "  -5.34%  "
